$p = $ppt.ActivePresentation

# --- Edit 1: update the table style id on the table in slide 15 ---
$tableSlide = $p.Slides.Item(15)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{A455968A-FAE1-40E4-B6BB-1E9309CD5A88}", $true)
    }
}

# --- Edit 2: add a new bullet paragraph at the end of the body text on slide 18 ---
$bodySlide = $p.Slides.Item(18)
for ($i = 1; $i -le $bodySlide.Shapes.Count; $i++) {
    $shp = $bodySlide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $lastText = $tr.Paragraphs($tr.Paragraphs().Count).Text
        if ($lastText -eq "Model all people involved (politicians, inspectors, judges, not just defendants)") {
            $lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
            [void]$lastPara.InsertAfter("`rPotential bias at every step")
        }
    }
}
